$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be interpreted as text so numeric-looking strings
# (e.g. "1.005") are not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.471.66"
$ws.Range("E2").Value = "  +1.67%  "
$ws.Range("D3").Value = "1.867.29"
$ws.Range("E3").Value = "  +2.40%  "
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "315.16"
$ws.Range("E5").Value = "  +1.89%  "
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").Value = "0.4668"
$ws.Range("E7").Value = "  -0.45%  "
$ws.Range("D8").Value = "0.3728"
$ws.Range("E8").Value = "  +1.86%  "
$ws.Range("D9").Value = "0.07375"
$ws.Range("E9").Value = "  +1.85%  "
$ws.Range("D10").Value = "0.8883"
$ws.Range("E10").Value = "  +3.19%  "
$ws.Range("D11").Value = "0.07920"
$ws.Range("E11").Value = "  +5.07%  "
$ws.Range("D12").Value = "20.01"
$ws.Range("E12").Value = "  +0.77%  "
$ws.Range("D13").Value = "1.839.23"
$ws.Range("E13").Value = "  +0.37%  "
$ws.Range("D14").Value = "5.422"
$ws.Range("E14").Value = "  +1.73%  "
$ws.Range("D15").Value = "6.610"
$ws.Range("E15").Value = "  +2.14%  "
$ws.Range("D16").Value = "92.81"
$ws.Range("E16").Value = "  +0.83%  "
$ws.Range("E17").Value = "  -0.14%  "
$ws.Range("D18").Value = "0.000008914"
$ws.Range("E18").Value = "  +3.14%  "
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("E20").Value = "  +2.97%  "
$ws.Range("D21").Value = "27.499.46"
$ws.Range("E21").Value = "  +3.57%  "
$ws.Range("D22").Value = "5.162"
$ws.Range("E22").Value = "  +0.25%  "
$ws.Range("E23").Value = "  +0.60%  "
$ws.Range("D24").Value = "2.088.24"
$ws.Range("E24").Value = "  -0.55%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "1.898"
$ws.Range("E25").Value = "  +3.08%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "153.42"
$ws.Range("E26").Value = "  +1.22%  "
$ws.Range("D27").Value = "18.54"
$ws.Range("E27").Value = "  +2.02%  "
$ws.Range("E28").Value = "  +0.88%  "
$ws.Range("D29").Value = "5.178"
$ws.Range("E29").Value = "  +1.03%  "
$ws.Range("D30").Value = "117.02"
$ws.Range("E30").Value = "  +1.44%  "
$ws.Range("D31").Value = "0.08909"
$ws.Range("E31").Value = "  +0.30%  "
$ws.Range("D32").Value = "0.7617"
$ws.Range("E32").Value = "  +5.86%  "
$ws.Range("D33").Value = "3.028"
$ws.Range("E33").Value = "  +2.31%  "
$ws.Range("D34").Value = "1.171"
$ws.Range("E34").Value = "  +3.53%  "
$ws.Range("D35").Value = "4.498"
$ws.Range("E35").Value = "  +1.66%  "
$ws.Range("D36").Value = "2.659"
$ws.Range("E36").Value = "  +10.10%  "
$ws.Range("D37").Value = "0.01973"
$ws.Range("E37").Value = "  +2.55%  "
$ws.Range("E38").Value = "  +0.32%  "
$ws.Range("D39").Value = "0.05277"
$ws.Range("E39").Value = "  +0.34%  "
$ws.Range("D40").Value = "2.999"
$ws.Range("E40").Value = "  +2.32%  "
$ws.Range("D41").Value = "7.190"
$ws.Range("E41").Value = "  +0.72%  "
$ws.Range("D42").Value = "0.5196"
$ws.Range("E42").Value = "  +0.65%  "
$ws.Range("E43").Value = "  +1.33%  "
$ws.Range("D44").Value = "8.387"
$ws.Range("E44").Value = "  +2.49%  "
$ws.Range("D45").Value = "0.4875"
$ws.Range("E45").Value = "  +1.08%  "
$ws.Range("D46").Value = "10.40"
$ws.Range("E46").Value = "  +2.72%  "
$ws.Range("E47").Value = "  -0.20%  "
$ws.Range("D48").Value = "104.25"
$ws.Range("E48").Value = "  +1.45%  "
$ws.Range("D49").Value = "1.658"
$ws.Range("E49").Value = "  +2.14%  "
$ws.Range("D50").Value = "0.06279"
$ws.Range("E50").Value = "  +0.50%  "
$ws.Range("E51").Value = "  +3.37%  "

# Restore default (General/Normal) styling on column D so no stray
# number-format style lingers on the cells.
$ws.Range("D2:D51").Style = "Normal"

